$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "1.2 Definicion del Problema" was demoted from Heading 3 to Heading 2
#    (pStyle Ttulo3 -> Ttulo2) for the heading that carries bookmark
#    _Toc444128961.
# ---------------------------------------------------------------------------
$problemHeading = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Heading 3" -and $p.Range.Text -like "*Definici*n del Problema*") {
        $problemHeading = $p
        break
    }
}
if ($problemHeading -ne $null) {
    $problemHeading.Style = $d.Styles("Heading 2")
}

# ---------------------------------------------------------------------------
# 2) The "_GoBack" bookmark (Word's last-edit-location marker) moved from
#    the end of the "...Torres, 2011)." paragraph to the end of the first
#    bullet under the "1.6 Alcance" heading ("...funcione via web.").
#    Removing it here and re-adding it there also shifts every bookmark id
#    that sits in between down by one (matching the diff's id 9 -> 8 for
#    _Toc444128967 and its neighbours).
# ---------------------------------------------------------------------------
try {
    $oldGoBack = $d.Bookmarks("_GoBack")
    $oldGoBack.Delete()
} catch {
    # no pre-existing _GoBack bookmark - nothing to remove
}

$scopeTarget = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*funcione v*a web*") {
        $scopeTarget = $p
        break
    }
}
if ($scopeTarget -ne $null) {
    # Insert a temporary marker right before the paragraph mark, wrap a
    # bookmark around it (collapsed ranges landing exactly on a paragraph
    # boundary are unreliable for Bookmarks.Add), then clear the marker
    # text again so the bookmark collapses to the correct, empty position.
    $insertPos = $scopeTarget.Range.End - 1
    $insertionPoint = $d.Range($insertPos, $insertPos)
    $insertionPoint.InsertAfter("TMPGOBACKMARKER")
    $markerRange = $d.Range($insertPos, $insertPos + 15)
    $markerRange.Bookmarks.Add("_GoBack") | Out-Null
    $markerRange.Text = ""
}

# ---------------------------------------------------------------------------
# 3) The footer's cached PAGE field result changed from "2" to "3".
# ---------------------------------------------------------------------------
$section = $d.Sections.Item(1)
$footer = $section.Footers.Item(1)
$footerRange = $footer.Range
for ($i = 1; $i -le $footerRange.Characters.Count; $i++) {
    $ch = $footerRange.Characters.Item($i)
    if ($ch.Text -eq "2") {
        $ch.Text = "3"
        break
    }
}

Write-Output "capituloI.docx subtitle corrections applied"
